$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 1.027
$ws.Cells.Item(2, 5).Value = 5.3
$ws.Cells.Item(2, 6).Value = 1.23625
$ws.Cells.Item(2, 7).Value = 0.09761414583333339
$ws.Cells.Item(2, 8).Value = 3.852075156483195
$ws.Cells.Item(2, 9).Value = 1.333033786264961
$ws.Cells.Item(2, 10).Value = 1.020353483153802
$ws.Cells.Item(2, 11).Value = 0.1573903470329096
$ws.Cells.Item(2, 12).Value = 0.08029689440490745

$ws.Cells.Item(3, 4).Value = 1.239
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(3, 6).Value = 1.23625
$ws.Cells.Item(3, 7).Value = 0.09761414583333339
$ws.Cells.Item(3, 8).Value = 3.678867773382904
$ws.Cells.Item(3, 9).Value = 1.333033786264961
$ws.Cells.Item(3, 10).Value = 1.27684112895035
$ws.Cells.Item(3, 11).Value = 0.3324226319460587
$ws.Cells.Item(3, 12).Value = 0.2122254443313261

$ws.Cells.Item(4, 4).Value = 1.166
$ws.Cells.Item(4, 5).Value = 2.5
$ws.Cells.Item(4, 6).Value = 1.23625
$ws.Cells.Item(4, 7).Value = 0.09761414583333339
$ws.Cells.Item(4, 8).Value = 3.738509938318381
$ws.Cells.Item(4, 9).Value = 1.333033786264961
$ws.Cells.Item(4, 10).Value = 1.245017365128281
$ws.Cells.Item(4, 11).Value = 0.1943661439160063
$ws.Cells.Item(4, 12).Value = 0.1209946121842252

$ws.Cells.Item(5, 4).Value = 0.9089999999999999
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1.23625
$ws.Cells.Item(5, 7).Value = 0.09761414583333339
$ws.Cells.Item(5, 8).Value = 3.948483039529583
$ws.Cells.Item(5, 9).Value = 1.333033786264961
$ws.Cells.Item(5, 10).Value = 0.7377704624374438
$ws.Cells.Item(5, 11).Value = 0.2465721666239024
$ws.Cells.Item(5, 12).Value = 0.09095683069715947

$ws.Cells.Item(6, 4).Value = 1.1865
$ws.Cells.Item(6, 5).Value = 2.5
$ws.Cells.Item(6, 6).Value = 1.23625
$ws.Cells.Item(6, 7).Value = 0.09761414583333339
$ws.Cells.Item(6, 8).Value = 3.721761111178966
$ws.Cells.Item(6, 9).Value = 1.333033786264961
$ws.Cells.Item(6, 10).Value = 1.260804627539021
$ws.Cells.Item(6, 11).Value = 0.1973935966133143
$ws.Cells.Item(6, 12).Value = 0.1244373800283187

$ws.Cells.Item(7, 4).Value = 1.89
$ws.Cells.Item(7, 5).Value = 1.8
$ws.Cells.Item(7, 6).Value = 1.23625
$ws.Cells.Item(7, 7).Value = 0.09761414583333339
$ws.Cells.Item(7, 8).Value = 3.146990384711729
$ws.Cells.Item(7, 9).Value = 1.333033786264961
$ws.Cells.Item(7, 10).Value = 0.143023266457987
$ws.Cells.Item(7, 11).Value = 0.1749572221769412
$ws.Cells.Item(7, 12).Value = 0.01251147670308094

$ws.Cells.Item(8, 4).Value = 2.171
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 1.694111111111111
$ws.Cells.Item(8, 7).Value = 0.2561913765432099
$ws.Cells.Item(8, 8).Value = 2.917408900508042
$ws.Cells.Item(8, 9).Value = 1.333033786264961
$ws.Cells.Item(8, 10).Value = 0.5056669539171873
$ws.Cells.Item(8, 11).Value = 0.06791642735171953
$ws.Cells.Item(8, 12).Value = 0.01717154646994098

$ws.Cells.Item(9, 4).Value = 1.8415
$ws.Cells.Item(9, 5).Value = 2.9
$ws.Cells.Item(9, 6).Value = 1.694111111111111
$ws.Cells.Item(9, 7).Value = 0.2561913765432099
$ws.Cells.Item(9, 8).Value = 3.186615658675711
$ws.Cells.Item(9, 9).Value = 1.333033786264961
$ws.Cells.Item(9, 10).Value = 0.7554661938984805
$ws.Cells.Item(9, 11).Value = 0.3350485440787231
$ws.Cells.Item(9, 12).Value = 0.1265589241831901

$ws.Cells.Item(10, 4).Value = 1.031
$ws.Cells.Item(10, 5).Value = 2.95
$ws.Cells.Item(10, 6).Value = 1.694111111111111
$ws.Cells.Item(10, 7).Value = 0.2561913765432099
$ws.Cells.Item(10, 8).Value = 3.848807092651113
$ws.Cells.Item(10, 9).Value = 1.333033786264961
$ws.Cells.Item(10, 10).Value = 0.3341374974572674
$ws.Cells.Item(10, 11).Value = 0.255206917186555
$ws.Cells.Item(10, 12).Value = 0.04263710032124978

$ws.Cells.Item(11, 4).Value = 1.742
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1.694111111111111
$ws.Cells.Item(11, 7).Value = 0.2561913765432099
$ws.Cells.Item(11, 8).Value = 3.267908746498724
$ws.Cells.Item(11, 9).Value = 1.333033786264961
$ws.Cells.Item(11, 10).Value = 0.7846644510953596
$ws.Cells.Item(11, 11).Value = 0.3363547285902926
$ws.Cells.Item(11, 12).Value = 0.1319627992413153

$ws.Cells.Item(12, 4).Value = 2.302
$ws.Cells.Item(12, 5).Value = 1.7
$ws.Cells.Item(12, 6).Value = 1.694111111111111
$ws.Cells.Item(12, 7).Value = 0.2561913765432099
$ws.Cells.Item(12, 8).Value = 2.810379810007391
$ws.Cells.Item(12, 9).Value = 1.333033786264961
$ws.Cells.Item(12, 10).Value = 0.3831915621817934
$ws.Cells.Item(12, 11).Value = 0.2175938845554888
$ws.Cells.Item(12, 12).Value = 0.04169007027201128

$ws.Cells.Item(13, 4).Value = 1.2415
$ws.Cells.Item(13, 5).Value = 4.9
$ws.Cells.Item(13, 6).Value = 1.694111111111111
$ws.Cells.Item(13, 7).Value = 0.2561913765432099
$ws.Cells.Item(13, 8).Value = 3.676825233487854
$ws.Cells.Item(13, 9).Value = 1.333033786264961
$ws.Cells.Item(13, 10).Value = 0.5284350661159578
$ws.Cells.Item(13, 11).Value = 0.1971378608669609
$ws.Cells.Item(13, 12).Value = 0.05208727927059548

$ws.Cells.Item(14, 4).Value = 1.07
$ws.Cells.Item(14, 5).Value = 4.25
$ws.Cells.Item(14, 6).Value = 1.694111111111111
$ws.Cells.Item(14, 7).Value = 0.2561913765432099
$ws.Cells.Item(14, 8).Value = 3.816943470288324
$ws.Cells.Item(14, 9).Value = 1.333033786264961
$ws.Cells.Item(14, 10).Value = 0.3685327064981622
$ws.Cells.Item(14, 11).Value = 0.3220624641982008
$ws.Cells.Item(14, 12).Value = 0.0593452757962152

$ws.Cells.Item(15, 4).Value = 2.446
$ws.Cells.Item(15, 5).Value = 2.55
$ws.Cells.Item(15, 6).Value = 1.694111111111111
$ws.Cells.Item(15, 7).Value = 0.2561913765432099
$ws.Cells.Item(15, 8).Value = 2.692729512052477
$ws.Cells.Item(15, 9).Value = 1.333033786264961
$ws.Cells.Item(15, 10).Value = 0.2614866445343001
$ws.Cells.Item(15, 11).Value = 0.3429027718760512
$ws.Cells.Item(15, 12).Value = 0.0448322476096896

$ws.Cells.Item(16, 4).Value = 1.402
$ws.Cells.Item(16, 5).Value = 5.5
$ws.Cells.Item(16, 6).Value = 1.694111111111111
$ws.Cells.Item(16, 7).Value = 0.2561913765432099
$ws.Cells.Item(16, 8).Value = 3.545694172225605
$ws.Cells.Item(16, 9).Value = 1.333033786264961
$ws.Cells.Item(16, 10).Value = 0.6672724976589723
$ws.Cells.Item(16, 11).Value = 0.08247741402644236
$ws.Cells.Item(16, 12).Value = 0.02751745502893868

Write-Host "Updated likelihood formula values for rows 2-16"